$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 1. Move the small picture ("Image 95") inside the "Groupe 37" group. ---
# Target raw <a:off> (child-space, i.e. same space as the group's chOff/chExt)
# is x=4871722 y=691894 (EMU). This interop's Shape.Left/Top setters for a
# grouped shape write the value (in points) straight into the raw <a:off>
# without re-deriving it through the group's child<->slide transform, so we
# simply convert the target EMU values to points (1 pt = 12700 EMU).
$grp = $s.Shapes.Item("Groupe 37")
$pic = $grp.GroupItems.Item("Image 95")
$pic.Left = 4871722 / 12700
$pic.Top = 691894 / 12700

# --- 2. Merge the "and" + " Time " runs into a single "and Time " run. ---
$shText = $s.Shapes.Item("Rectangle à coins arrondis 17")
$tr = $shText.TextFrame.TextRange
$paraIdx = 0
foreach ($para in $tr.Paragraphs()) {
    $paraIdx++
    if ($paraIdx -eq 2) {
        $runIdx = 0
        foreach ($run in $para.Runs()) {
            $runIdx++
            if ($runIdx -eq 1) {
                $run.Text = "and Time "
            } elseif ($runIdx -eq 2) {
                $run.Text = ""
            }
        }
    }
}
